$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# update 72 µm to 74 µm in the "16 to >72" / "17 to >72" range labels
$ws.Range("I5").Value = "16 to >74"
$ws.Range("I3").Value = "17 to >74"

# match the author's final selection on the sheet
$ws.Range("I3").Select()
